# BOT; UPDATE DATA
# Refreshes 2020-04-10 (serial 43949) PCR/infection-survey figures across
# the "all", "kobe" and "other" sheets, updates the "under investigation"
# caption on "all", and leaves the selection where Excel would after the
# edit (the last cell touched on each sheet).

$wb = $excel.ActiveWorkbook

# --- Sheet "all" ----------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")

# "under investigation" count caption
$wsAll.Range("B23").Value = "※14件調査中"

# Row 21 (2020-04-10) revised totals
$wsAll.Range("C21").Value = 241
$wsAll.Range("D21").Value = 129
$wsAll.Range("E21").Value = 120
$wsAll.Range("F21").Value = 9
$wsAll.Range("G21").Value = 4
$wsAll.Range("H21").Value = 108

$wsAll.Activate()
$wsAll.Range("A21").Select()

# --- Sheet "kobe" -----------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")

# Row 76 (2020-04-10) revised totals
$wsKobe.Range("B76").Value = 82
$wsKobe.Range("C76").Value = 1854
$wsKobe.Range("F76").Value = 125
$wsKobe.Range("G76").Value = 117
$wsKobe.Range("H76").Value = 8
$wsKobe.Range("I76").Value = 4
$wsKobe.Range("J76").Value = 101

$wsKobe.Activate()
$wsKobe.Range("A76").Select()
# Scroll the frozen (bottom-right) pane so row 58 sits at the top of view.
$excel.ActiveWindow.ScrollRow = 58
$excel.ActiveWindow.ScrollColumn = 2

# --- Sheet "other" ----------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")

# Row 51 (2020-04-10) revised totals
$wsOther.Range("D51").Value = 4
$wsOther.Range("E51").Value = 3
$wsOther.Range("H51").Value = 7

$wsOther.Activate()
$wsOther.Range("A51").Select()
# Scroll the frozen (bottom-right) pane so column B sits at the left of view.
$excel.ActiveWindow.ScrollRow = 39
$excel.ActiveWindow.ScrollColumn = 2

# Leave "all" as the active sheet, matching the workbook's tabSelected state.
$wsAll.Activate()
